# Update the "取得日時" (acquired timestamp) column for all data rows
# on the "ランサーズ" sheet from "2025-10-29 12:39:32" to "2025-10-29 12:50:21".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-10-29 12:39:32"
$newValue = "2025-10-29 12:50:21"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
